$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp label (A1)
$ws.Range("A1").Value = 'Datos actualizados a 22 de Marzo de 2020 a las 17:46'

$data = New-Object 'object[,]' 189,8
$data[0,0] = 'China'
$data[0,1] = 81054
$data[0,2] = 46
$data[0,3] = 72440
$data[0,4] = 5353
$data[0,5] = 1845
$data[0,6] = 6
$data[0,7] = 3261
$data[1,0] = 'Italia'
$data[1,1] = 53578
$data[1,2] = 0
$data[1,3] = 6072
$data[1,4] = 42681
$data[1,5] = 2857
$data[1,6] = 0
$data[1,7] = 4825
$data[2,0] = 'Estados Unidos'
$data[2,1] = 30239
$data[2,2] = 6032
$data[2,3] = 178
$data[2,4] = 29673
$data[2,5] = 708
$data[2,6] = 86
$data[2,7] = 388
$data[3,0] = 'España'
$data[3,1] = 28603
$data[3,2] = 3107
$data[3,3] = 2125
$data[3,4] = 24722
$data[3,5] = 1785
$data[3,6] = 375
$data[3,7] = 1756
$data[4,0] = 'Alemania'
$data[4,1] = 23974
$data[4,2] = 1610
$data[4,3] = 266
$data[4,4] = 23615
$data[4,5] = 2
$data[4,6] = 9
$data[4,7] = 93
$data[5,0] = 'Iran'
$data[5,1] = 21638
$data[5,2] = 1028
$data[5,3] = 7635
$data[5,4] = 12318
$data[5,5] = 0
$data[5,6] = 129
$data[5,7] = 1685
$data[6,0] = 'Francia'
$data[6,1] = 14459
$data[6,2] = 0
$data[6,3] = 1587
$data[6,4] = 12310
$data[6,5] = 1525
$data[6,6] = 0
$data[6,7] = 562
$data[7,0] = 'Corea del Sur'
$data[7,1] = 8897
$data[7,2] = 98
$data[7,3] = 2909
$data[7,4] = 5884
$data[7,5] = 59
$data[7,6] = 2
$data[7,7] = 104
$data[8,0] = 'Suiza'
$data[8,1] = 7230
$data[8,2] = 367
$data[8,3] = 131
$data[8,4] = 7014
$data[8,5] = 141
$data[8,6] = 5
$data[8,7] = 85
$data[9,0] = 'Reino Unido'
$data[9,1] = 5018
$data[9,2] = 0
$data[9,3] = 93
$data[9,4] = 4681
$data[9,5] = 20
$data[9,6] = 11
$data[9,7] = 244
$data[10,0] = 'Paises Bajos'
$data[10,1] = 4204
$data[10,2] = 573
$data[10,3] = 2
$data[10,4] = 4023
$data[10,5] = 354
$data[10,6] = 43
$data[10,7] = 179
$data[11,0] = 'Belgica'
$data[11,1] = 3401
$data[11,2] = 586
$data[11,3] = 263
$data[11,4] = 3063
$data[11,5] = 288
$data[11,6] = 8
$data[11,7] = 75
$data[12,0] = 'Austria'
$data[12,1] = 3302
$data[12,2] = 310
$data[12,3] = 9
$data[12,4] = 3277
$data[12,5] = 15
$data[12,6] = 8
$data[12,7] = 16
$data[13,0] = 'Noruega'
$data[13,1] = 2262
$data[13,2] = 98
$data[13,3] = 6
$data[13,4] = 2249
$data[13,5] = 28
$data[13,6] = 0
$data[13,7] = 7
$data[14,0] = 'Suecia'
$data[14,1] = 1931
$data[14,2] = 161
$data[14,3] = 16
$data[14,4] = 1894
$data[14,5] = 68
$data[14,6] = 1
$data[14,7] = 21
$data[15,0] = 'Portugal'
$data[15,1] = 1600
$data[15,2] = 320
$data[15,3] = 5
$data[15,4] = 1581
$data[15,5] = 26
$data[15,6] = 2
$data[15,7] = 14
$data[16,0] = 'Dinamarca'
$data[16,1] = 1395
$data[16,2] = 69
$data[16,3] = 1
$data[16,4] = 1381
$data[16,5] = 42
$data[16,6] = 0
$data[16,7] = 13
$data[17,0] = 'Canada'
$data[17,1] = 1385
$data[17,2] = 57
$data[17,3] = 14
$data[17,4] = 1351
$data[17,5] = 1
$data[17,6] = 1
$data[17,7] = 20
$data[18,0] = 'Australia'
$data[18,1] = 1353
$data[18,2] = 281
$data[18,3] = 46
$data[18,4] = 1300
$data[18,5] = 2
$data[18,6] = 0
$data[18,7] = 7
$data[19,0] = 'Malasia'
$data[19,1] = 1306
$data[19,2] = 123
$data[19,3] = 139
$data[19,4] = 1157
$data[19,5] = 26
$data[19,6] = 2
$data[19,7] = 10
$data[20,0] = 'Brasil'
$data[20,1] = 1209
$data[20,2] = 31
$data[20,3] = 2
$data[20,4] = 1189
$data[20,5] = 18
$data[20,6] = 0
$data[20,7] = 18
$data[21,0] = 'Japon'
$data[21,1] = 1086
$data[21,2] = 32
$data[21,3] = 235
$data[21,4] = 815
$data[21,5] = 57
$data[21,6] = 0
$data[21,7] = 36
$data[22,0] = 'Chequia'
$data[22,1] = 1047
$data[22,2] = 52
$data[22,3] = 6
$data[22,4] = 1041
$data[22,5] = 19
$data[22,6] = 0
$data[22,7] = 0
$data[23,0] = 'Turquia'
$data[23,1] = 947
$data[23,2] = 0
$data[23,3] = 0
$data[23,4] = 926
$data[23,5] = 0
$data[23,6] = 0
$data[23,7] = 21
$data[24,0] = 'Israel'
$data[24,1] = 945
$data[24,2] = 62
$data[24,3] = 37
$data[24,4] = 907
$data[24,5] = 20
$data[24,6] = 0
$data[24,7] = 1
$data[25,0] = 'Luxemburgo'
$data[25,1] = 798
$data[25,2] = 128
$data[25,3] = 6
$data[25,4] = 784
$data[25,5] = 3
$data[25,6] = 0
$data[25,7] = 8
$data[26,0] = 'Ecuador'
$data[26,1] = 789
$data[26,2] = 257
$data[26,3] = 3
$data[26,4] = 772
$data[26,5] = 2
$data[26,6] = 7
$data[26,7] = 14
$data[27,0] = 'Irlanda'
$data[27,1] = 785
$data[27,2] = 0
$data[27,3] = 5
$data[27,4] = 777
$data[27,5] = 13
$data[27,6] = 0
$data[27,7] = 3
$data[28,0] = 'Crucero'
$data[28,1] = 712
$data[28,2] = 0
$data[28,3] = 567
$data[28,4] = 137
$data[28,5] = 15
$data[28,6] = 0
$data[28,7] = 8
$data[29,0] = 'Pakistan'
$data[29,1] = 646
$data[29,2] = 1
$data[29,3] = 13
$data[29,4] = 630
$data[29,5] = 0
$data[29,6] = 0
$data[29,7] = 3
$data[30,0] = 'Chile'
$data[30,1] = 632
$data[30,2] = 95
$data[30,3] = 8
$data[30,4] = 623
$data[30,5] = 7
$data[30,6] = 0
$data[30,7] = 1
$data[31,0] = 'Finlandia'
$data[31,1] = 626
$data[31,2] = 103
$data[31,3] = 10
$data[31,4] = 615
$data[31,5] = 12
$data[31,6] = 0
$data[31,7] = 1
$data[32,0] = 'Grecia'
$data[32,1] = 624
$data[32,2] = 94
$data[32,3] = 19
$data[32,4] = 590
$data[32,5] = 18
$data[32,6] = 2
$data[32,7] = 15
$data[33,0] = 'Tailandia'
$data[33,1] = 599
$data[33,2] = 188
$data[33,3] = 44
$data[33,4] = 554
$data[33,5] = 7
$data[33,6] = 0
$data[33,7] = 1
$data[34,0] = 'Islandia'
$data[34,1] = 568
$data[34,2] = 95
$data[34,3] = 5
$data[34,4] = 562
$data[34,5] = 1
$data[34,6] = 0
$data[34,7] = 1
$data[35,0] = 'Polonia'
$data[35,1] = 563
$data[35,2] = 27
$data[35,3] = 13
$data[35,4] = 543
$data[35,5] = 3
$data[35,6] = 2
$data[35,7] = 7
$data[36,0] = 'Indonesia'
$data[36,1] = 514
$data[36,2] = 64
$data[36,3] = 29
$data[36,4] = 437
$data[36,5] = 0
$data[36,6] = 10
$data[36,7] = 48
$data[37,0] = 'Arabia Saudita'
$data[37,1] = 511
$data[37,2] = 119
$data[37,3] = 17
$data[37,4] = 494
$data[37,5] = 0
$data[37,6] = 0
$data[37,7] = 0
$data[38,0] = 'Catar'
$data[38,1] = 481
$data[38,2] = 0
$data[38,3] = 27
$data[38,4] = 454
$data[38,5] = 6
$data[38,6] = 0
$data[38,7] = 0
$data[39,0] = 'Singapur'
$data[39,1] = 455
$data[39,2] = 23
$data[39,3] = 144
$data[39,4] = 309
$data[39,5] = 14
$data[39,6] = 0
$data[39,7] = 2
$data[40,0] = 'Rumania'
$data[40,1] = 433
$data[40,2] = 66
$data[40,3] = 64
$data[40,4] = 367
$data[40,5] = 14
$data[40,6] = 2
$data[40,7] = 2
$data[41,0] = 'Eslovenia'
$data[41,1] = 414
$data[41,2] = 31
$data[41,3] = 0
$data[41,4] = 412
$data[41,5] = 12
$data[41,6] = 1
$data[41,7] = 2
$data[42,0] = 'India'
$data[42,1] = 391
$data[42,2] = 59
$data[42,3] = 24
$data[42,4] = 360
$data[42,5] = 0
$data[42,6] = 2
$data[42,7] = 7
$data[43,0] = 'Filipinas'
$data[43,1] = 380
$data[43,2] = 73
$data[43,3] = 15
$data[43,4] = 340
$data[43,5] = 1
$data[43,6] = 6
$data[43,7] = 25
$data[44,0] = 'Rusia'
$data[44,1] = 367
$data[44,2] = 61
$data[44,3] = 16
$data[44,4] = 350
$data[44,5] = 0
$data[44,6] = 0
$data[44,7] = 1
$data[45,0] = 'Barein'
$data[45,1] = 332
$data[45,2] = 22
$data[45,3] = 149
$data[45,4] = 181
$data[45,5] = 4
$data[45,6] = 1
$data[45,7] = 2
$data[46,0] = 'Estonia'
$data[46,1] = 326
$data[46,2] = 20
$data[46,3] = 4
$data[46,4] = 322
$data[46,5] = 2
$data[46,6] = 0
$data[46,7] = 0
$data[47,0] = 'Peru'
$data[47,1] = 318
$data[47,2] = 0
$data[47,3] = 1
$data[47,4] = 312
$data[47,5] = 5
$data[47,6] = 0
$data[47,7] = 5
$data[48,0] = 'Hong Kong'
$data[48,1] = 317
$data[48,2] = 43
$data[48,3] = 100
$data[48,4] = 213
$data[48,5] = 4
$data[48,6] = 0
$data[48,7] = 4
$data[49,0] = 'Egipto'
$data[49,1] = 294
$data[49,2] = 0
$data[49,3] = 42
$data[49,4] = 242
$data[49,5] = 0
$data[49,6] = 0
$data[49,7] = 10
$data[50,0] = 'Croacia'
$data[50,1] = 254
$data[50,2] = 48
$data[50,3] = 5
$data[50,4] = 248
$data[50,5] = 5
$data[50,6] = 0
$data[50,7] = 1
$data[51,0] = 'Mexico'
$data[51,1] = 251
$data[51,2] = 48
$data[51,3] = 4
$data[51,4] = 245
$data[51,5] = 1
$data[51,6] = 0
$data[51,7] = 2
$data[52,0] = 'Libano'
$data[52,1] = 248
$data[52,2] = 18
$data[52,3] = 8
$data[52,4] = 236
$data[52,5] = 4
$data[52,6] = 0
$data[52,7] = 4
$data[53,0] = 'Panama'
$data[53,1] = 245
$data[53,2] = 0
$data[53,3] = 1
$data[53,4] = 241
$data[53,5] = 7
$data[53,6] = 0
$data[53,7] = 3
$data[54,0] = 'Sudafrica'
$data[54,1] = 240
$data[54,2] = 0
$data[54,3] = 2
$data[54,4] = 238
$data[54,5] = 0
$data[54,6] = 0
$data[54,7] = 0
$data[55,0] = 'Irak'
$data[55,1] = 233
$data[55,2] = 19
$data[55,3] = 57
$data[55,4] = 156
$data[55,5] = 0
$data[55,6] = 3
$data[55,7] = 20
$data[56,0] = 'Colombia'
$data[56,1] = 231
$data[56,2] = 35
$data[56,3] = 3
$data[56,4] = 226
$data[56,5] = 0
$data[56,6] = 2
$data[56,7] = 2
$data[57,0] = 'Argentina'
$data[57,1] = 225
$data[57,2] = 67
$data[57,3] = 27
$data[57,4] = 194
$data[57,5] = 0
$data[57,6] = 0
$data[57,7] = 4
$data[58,0] = 'Republica Dominicana'
$data[58,1] = 202
$data[58,2] = 90
$data[58,3] = 0
$data[58,4] = 199
$data[58,5] = 0
$data[58,6] = 0
$data[58,7] = 3
$data[59,0] = 'Argelia'
$data[59,1] = 201
$data[59,2] = 62
$data[59,3] = 65
$data[59,4] = 119
$data[59,5] = 0
$data[59,6] = 2
$data[59,7] = 17
$data[60,0] = 'Armenia'
$data[60,1] = 190
$data[60,2] = 30
$data[60,3] = 2
$data[60,4] = 188
$data[60,5] = 6
$data[60,6] = 0
$data[60,7] = 0
$data[61,0] = 'Serbia'
$data[61,1] = 188
$data[61,2] = 17
$data[61,3] = 2
$data[61,4] = 184
$data[61,5] = 4
$data[61,6] = 1
$data[61,7] = 2
$data[62,0] = 'Kuwait'
$data[62,1] = 188
$data[62,2] = 12
$data[62,3] = 27
$data[62,4] = 161
$data[62,5] = 5
$data[62,6] = 0
$data[62,7] = 0
$data[63,0] = 'Bulgaria'
$data[63,1] = 185
$data[63,2] = 22
$data[63,3] = 3
$data[63,4] = 179
$data[63,5] = 3
$data[63,6] = 0
$data[63,7] = 3
$data[64,0] = 'Eslovaquia'
$data[64,1] = 178
$data[64,2] = 0
$data[64,3] = 7
$data[64,4] = 171
$data[64,5] = 2
$data[64,6] = 0
$data[64,7] = 0
$data[65,0] = 'San Marino'
$data[65,1] = 175
$data[65,2] = 15
$data[65,3] = 4
$data[65,4] = 151
$data[65,5] = 13
$data[65,6] = 0
$data[65,7] = 20
$data[66,0] = 'Taiwan'
$data[66,1] = 169
$data[66,2] = 16
$data[66,3] = 28
$data[66,4] = 139
$data[66,5] = 0
$data[66,6] = 0
$data[66,7] = 2
$data[67,0] = 'Emiratos Arabes Unidos'
$data[67,1] = 153
$data[67,2] = 0
$data[67,3] = 38
$data[67,4] = 113
$data[67,5] = 2
$data[67,6] = 0
$data[67,7] = 2
$data[68,0] = 'Letonia'
$data[68,1] = 139
$data[68,2] = 15
$data[68,3] = 1
$data[68,4] = 138
$data[68,5] = 0
$data[68,6] = 0
$data[68,7] = 0
$data[69,0] = 'Uruguay'
$data[69,1] = 135
$data[69,2] = 25
$data[69,3] = 0
$data[69,4] = 135
$data[69,5] = 2
$data[69,6] = 0
$data[69,7] = 0
$data[70,0] = 'Hungria'
$data[70,1] = 131
$data[70,2] = 28
$data[70,3] = 16
$data[70,4] = 109
$data[70,5] = 6
$data[70,6] = 2
$data[70,7] = 6
$data[71,0] = 'Lituania'
$data[71,1] = 129
$data[71,2] = 30
$data[71,3] = 1
$data[71,4] = 127
$data[71,5] = 1
$data[71,6] = 0
$data[71,7] = 1
$data[72,0] = 'Costa Rica'
$data[72,1] = 117
$data[72,2] = 0
$data[72,3] = 2
$data[72,4] = 113
$data[72,5] = 2
$data[72,6] = 0
$data[72,7] = 2
$data[73,0] = 'Republica de Macedonia'
$data[73,1] = 115
$data[73,2] = 30
$data[73,3] = 1
$data[73,4] = 113
$data[73,5] = 1
$data[73,6] = 1
$data[73,7] = 1
$data[74,0] = 'Islas Feroe'
$data[74,1] = 115
$data[74,2] = 23
$data[74,3] = 3
$data[74,4] = 112
$data[74,5] = 0
$data[74,6] = 0
$data[74,7] = 0
$data[75,0] = 'Principado de Andorra'
$data[75,1] = 113
$data[75,2] = 25
$data[75,3] = 1
$data[75,4] = 111
$data[75,5] = 2
$data[75,6] = 1
$data[75,7] = 1
$data[76,0] = 'Vietnam'
$data[76,1] = 113
$data[76,2] = 19
$data[76,3] = 17
$data[76,4] = 96
$data[76,5] = 2
$data[76,6] = 0
$data[76,7] = 0
$data[77,0] = 'Marruecos'
$data[77,1] = 109
$data[77,2] = 13
$data[77,3] = 3
$data[77,4] = 103
$data[77,5] = 1
$data[77,6] = 0
$data[77,7] = 3
$data[78,0] = 'Jordania'
$data[78,1] = 100
$data[78,2] = 0
$data[78,3] = 1
$data[78,4] = 99
$data[78,5] = 0
$data[78,6] = 0
$data[78,7] = 0
$data[79,0] = 'Republica de Chipre'
$data[79,1] = 95
$data[79,2] = 11
$data[79,3] = 3
$data[79,4] = 91
$data[79,5] = 3
$data[79,6] = 0
$data[79,7] = 1
$data[80,0] = 'Moldavia'
$data[80,1] = 94
$data[80,2] = 14
$data[80,3] = 1
$data[80,4] = 92
$data[80,5] = 3
$data[80,6] = 0
$data[80,7] = 1
$data[81,0] = 'Bosnia y Herzegovina'
$data[81,1] = 94
$data[81,2] = 1
$data[81,3] = 2
$data[81,4] = 91
$data[81,5] = 1
$data[81,6] = 0
$data[81,7] = 1
$data[82,0] = 'Malta'
$data[82,1] = 90
$data[82,2] = 17
$data[82,3] = 2
$data[82,4] = 88
$data[82,5] = 1
$data[82,6] = 0
$data[82,7] = 0
$data[83,0] = 'Albania'
$data[83,1] = 89
$data[83,2] = 13
$data[83,3] = 2
$data[83,4] = 85
$data[83,5] = 2
$data[83,6] = 0
$data[83,7] = 2
$data[84,0] = 'Brunei'
$data[84,1] = 88
$data[84,2] = 5
$data[84,3] = 2
$data[84,4] = 86
$data[84,5] = 2
$data[84,6] = 0
$data[84,7] = 0
$data[85,0] = 'Camboya'
$data[85,1] = 84
$data[85,2] = 31
$data[85,3] = 2
$data[85,4] = 82
$data[85,5] = 0
$data[85,6] = 0
$data[85,7] = 0
$data[86,0] = 'Sri Lanka'
$data[86,1] = 82
$data[86,2] = 5
$data[86,3] = 3
$data[86,4] = 79
$data[86,5] = 2
$data[86,6] = 0
$data[86,7] = 0
$data[87,0] = 'Bielorrusia'
$data[87,1] = 76
$data[87,2] = 0
$data[87,3] = 15
$data[87,4] = 61
$data[87,5] = 0
$data[87,6] = 0
$data[87,7] = 0
$data[88,0] = 'Tunez'
$data[88,1] = 75
$data[88,2] = 15
$data[88,3] = 1
$data[88,4] = 71
$data[88,5] = 7
$data[88,6] = 2
$data[88,7] = 3
$data[89,0] = 'Burkina Faso'
$data[89,1] = 75
$data[89,2] = 11
$data[89,3] = 5
$data[89,4] = 66
$data[89,5] = 0
$data[89,6] = 1
$data[89,7] = 4
$data[90,0] = 'Venezuela'
$data[90,1] = 70
$data[90,2] = 0
$data[90,3] = 15
$data[90,4] = 55
$data[90,5] = 2
$data[90,6] = 0
$data[90,7] = 0
$data[91,0] = 'Nueva Zelanda'
$data[91,1] = 66
$data[91,2] = 14
$data[91,3] = 0
$data[91,4] = 66
$data[91,5] = 0
$data[91,6] = 0
$data[91,7] = 0
$data[92,0] = 'Azerbaiyan'
$data[92,1] = 65
$data[92,2] = 12
$data[92,3] = 11
$data[92,4] = 53
$data[92,5] = 0
$data[92,6] = 0
$data[92,7] = 1
$data[93,0] = 'Kazajistan'
$data[93,1] = 59
$data[93,2] = 5
$data[93,3] = 0
$data[93,4] = 59
$data[93,5] = 0
$data[93,6] = 0
$data[93,7] = 0
$data[94,0] = 'Estado de Palestina'
$data[94,1] = 59
$data[94,2] = 6
$data[94,3] = 17
$data[94,4] = 42
$data[94,5] = 0
$data[94,6] = 0
$data[94,7] = 0
$data[95,0] = 'Guadalupe'
$data[95,1] = 56
$data[95,2] = 0
$data[95,3] = 0
$data[95,4] = 55
$data[95,5] = 4
$data[95,6] = 0
$data[95,7] = 1
$data[96,0] = 'Senegal'
$data[96,1] = 56
$data[96,2] = 0
$data[96,3] = 5
$data[96,4] = 51
$data[96,5] = 0
$data[96,6] = 0
$data[96,7] = 0
$data[97,0] = 'Oman'
$data[97,1] = 55
$data[97,2] = 3
$data[97,3] = 17
$data[97,4] = 38
$data[97,5] = 0
$data[97,6] = 0
$data[97,7] = 0
$data[98,0] = 'Georgia'
$data[98,1] = 54
$data[98,2] = 5
$data[98,3] = 1
$data[98,4] = 53
$data[98,5] = 1
$data[98,6] = 0
$data[98,7] = 0
$data[99,0] = 'Trinidad yTobago'
$data[99,1] = 50
$data[99,2] = 1
$data[99,3] = 0
$data[99,4] = 50
$data[99,5] = 0
$data[99,6] = 0
$data[99,7] = 0
$data[100,0] = 'Reunion'
$data[100,1] = 47
$data[100,2] = 0
$data[100,3] = 1
$data[100,4] = 46
$data[100,5] = 0
$data[100,6] = 0
$data[100,7] = 0
$data[101,0] = 'Ucrania'
$data[101,1] = 47
$data[101,2] = 0
$data[101,3] = 1
$data[101,4] = 43
$data[101,5] = 0
$data[101,6] = 0
$data[101,7] = 3
$data[102,0] = 'Uzbekistan'
$data[102,1] = 43
$data[102,2] = 2
$data[102,3] = 0
$data[102,4] = 43
$data[102,5] = 0
$data[102,6] = 0
$data[102,7] = 0
$data[103,0] = 'Camerun'
$data[103,1] = 40
$data[103,2] = 0
$data[103,3] = 2
$data[103,4] = 38
$data[103,5] = 0
$data[103,6] = 0
$data[103,7] = 0
$data[104,0] = 'Liechtenstein'
$data[104,1] = 37
$data[104,2] = 0
$data[104,3] = 0
$data[104,4] = 37
$data[104,5] = 0
$data[104,6] = 0
$data[104,7] = 0
$data[105,0] = 'Martinica'
$data[105,1] = 37
$data[105,2] = 0
$data[105,3] = 0
$data[105,4] = 36
$data[105,5] = 7
$data[105,6] = 0
$data[105,7] = 1
$data[106,0] = 'Afganistan'
$data[106,1] = 34
$data[106,2] = 10
$data[106,3] = 1
$data[106,4] = 32
$data[106,5] = 0
$data[106,6] = 1
$data[106,7] = 1
$data[107,0] = 'Consejo Danes para los Refugiados'
$data[107,1] = 30
$data[107,2] = 7
$data[107,3] = 0
$data[107,4] = 29
$data[107,5] = 0
$data[107,6] = 0
$data[107,7] = 1
$data[108,0] = 'Guam'
$data[108,1] = 27
$data[108,2] = 12
$data[108,3] = 0
$data[108,4] = 26
$data[108,5] = 0
$data[108,6] = 1
$data[108,7] = 1
$data[109,0] = 'Nigeria'
$data[109,1] = 27
$data[109,2] = 5
$data[109,3] = 2
$data[109,4] = 25
$data[109,5] = 0
$data[109,6] = 0
$data[109,7] = 0
$data[110,0] = 'Banglades'
$data[110,1] = 27
$data[110,2] = 3
$data[110,3] = 3
$data[110,4] = 22
$data[110,5] = 0
$data[110,6] = 0
$data[110,7] = 2
$data[111,0] = 'Honduras'
$data[111,1] = 26
$data[111,2] = 2
$data[111,3] = 0
$data[111,4] = 26
$data[111,5] = 0
$data[111,6] = 0
$data[111,7] = 0
$data[112,0] = 'Bolivia'
$data[112,1] = 24
$data[112,2] = 5
$data[112,3] = 0
$data[112,4] = 24
$data[112,5] = 0
$data[112,6] = 0
$data[112,7] = 0
$data[113,0] = 'Mauricio'
$data[113,1] = 24
$data[113,2] = 10
$data[113,3] = 0
$data[113,4] = 22
$data[113,5] = 1
$data[113,6] = 1
$data[113,7] = 2
$data[114,0] = 'Puerto Rico'
$data[114,1] = 23
$data[114,2] = 2
$data[114,3] = 0
$data[114,4] = 22
$data[114,5] = 0
$data[114,6] = 0
$data[114,7] = 1
$data[115,0] = 'Paraguay'
$data[115,1] = 22
$data[115,2] = 0
$data[115,3] = 0
$data[115,4] = 21
$data[115,5] = 1
$data[115,6] = 0
$data[115,7] = 1
$data[116,0] = 'Macao'
$data[116,1] = 22
$data[116,2] = 3
$data[116,3] = 10
$data[116,4] = 12
$data[116,5] = 0
$data[116,6] = 0
$data[116,7] = 0
$data[117,0] = 'Cuba'
$data[117,1] = 21
$data[117,2] = 0
$data[117,3] = 0
$data[117,4] = 20
$data[117,5] = 0
$data[117,6] = 0
$data[117,7] = 1
$data[118,0] = 'Ghana'
$data[118,1] = 21
$data[118,2] = 0
$data[118,3] = 0
$data[118,4] = 20
$data[118,5] = 0
$data[118,6] = 0
$data[118,7] = 1
$data[119,0] = 'Jamaica'
$data[119,1] = 19
$data[119,2] = 0
$data[119,3] = 2
$data[119,4] = 16
$data[119,5] = 0
$data[119,6] = 0
$data[119,7] = 1
$data[120,0] = 'Guayana Francesa'
$data[120,1] = 18
$data[120,2] = 0
$data[120,3] = 0
$data[120,4] = 18
$data[120,5] = 0
$data[120,6] = 0
$data[120,7] = 0
$data[121,0] = 'Guyana'
$data[121,1] = 18
$data[121,2] = 0
$data[121,3] = 0
$data[121,4] = 17
$data[121,5] = 0
$data[121,6] = 0
$data[121,7] = 1
$data[122,0] = 'Monaco'
$data[122,1] = 18
$data[122,2] = 0
$data[122,3] = 1
$data[122,4] = 17
$data[122,5] = 0
$data[122,6] = 0
$data[122,7] = 0
$data[123,0] = 'Ruanda'
$data[123,1] = 17
$data[123,2] = 0
$data[123,3] = 0
$data[123,4] = 17
$data[123,5] = 0
$data[123,6] = 0
$data[123,7] = 0
$data[124,0] = 'Guatemala'
$data[124,1] = 17
$data[124,2] = 0
$data[124,3] = 0
$data[124,4] = 16
$data[124,5] = 0
$data[124,6] = 0
$data[124,7] = 1
$data[125,0] = 'Montenegro'
$data[125,1] = 16
$data[125,2] = 0
$data[125,3] = 0
$data[125,4] = 16
$data[125,5] = 0
$data[125,6] = 0
$data[125,7] = 0
$data[126,0] = 'Togo'
$data[126,1] = 16
$data[126,2] = 0
$data[126,3] = 0
$data[126,4] = 16
$data[126,5] = 0
$data[126,6] = 0
$data[126,7] = 0
$data[127,0] = 'Kenia'
$data[127,1] = 15
$data[127,2] = 8
$data[127,3] = 0
$data[127,4] = 15
$data[127,5] = 0
$data[127,6] = 0
$data[127,7] = 0
$data[128,0] = 'Polinesia Francesa'
$data[128,1] = 15
$data[128,2] = 0
$data[128,3] = 0
$data[128,4] = 15
$data[128,5] = 0
$data[128,6] = 0
$data[128,7] = 0
$data[129,0] = 'Gibraltar'
$data[129,1] = 15
$data[129,2] = 5
$data[129,3] = 2
$data[129,4] = 13
$data[129,5] = 0
$data[129,6] = 0
$data[129,7] = 0
$data[130,0] = 'Kirguistan'
$data[130,1] = 14
$data[130,2] = 0
$data[130,3] = 0
$data[130,4] = 14
$data[130,5] = 0
$data[130,6] = 0
$data[130,7] = 0
$data[131,0] = 'Barbados'
$data[131,1] = 14
$data[131,2] = 0
$data[131,3] = 0
$data[131,4] = 14
$data[131,5] = 0
$data[131,6] = 0
$data[131,7] = 0
$data[132,0] = 'Costa de Marfil'
$data[132,1] = 14
$data[132,2] = 0
$data[132,3] = 1
$data[132,4] = 13
$data[132,5] = 0
$data[132,6] = 0
$data[132,7] = 0
$data[133,0] = 'Maldivas'
$data[133,1] = 13
$data[133,2] = 0
$data[133,3] = 3
$data[133,4] = 10
$data[133,5] = 0
$data[133,6] = 0
$data[133,7] = 0
$data[134,0] = 'Tanzania'
$data[134,1] = 12
$data[134,2] = 6
$data[134,3] = 0
$data[134,4] = 12
$data[134,5] = 0
$data[134,6] = 0
$data[134,7] = 0
$data[135,0] = 'Etiopia'
$data[135,1] = 11
$data[135,2] = 2
$data[135,3] = 0
$data[135,4] = 11
$data[135,5] = 0
$data[135,6] = 0
$data[135,7] = 0
$data[136,0] = 'Mayotte'
$data[136,1] = 11
$data[136,2] = 0
$data[136,3] = 0
$data[136,4] = 11
$data[136,5] = 0
$data[136,6] = 0
$data[136,7] = 0
$data[137,0] = 'Mongolia'
$data[137,1] = 10
$data[137,2] = 0
$data[137,3] = 0
$data[137,4] = 10
$data[137,5] = 0
$data[137,6] = 0
$data[137,7] = 0
$data[138,0] = 'Aruba'
$data[138,1] = 8
$data[138,2] = 3
$data[138,3] = 1
$data[138,4] = 7
$data[138,5] = 0
$data[138,6] = 0
$data[138,7] = 0
$data[139,0] = 'Seychelles'
$data[139,1] = 7
$data[139,2] = 0
$data[139,3] = 0
$data[139,4] = 7
$data[139,5] = 0
$data[139,6] = 0
$data[139,7] = 0
$data[140,0] = 'Guinea Ecuatorial'
$data[140,1] = 6
$data[140,2] = 0
$data[140,3] = 0
$data[140,4] = 6
$data[140,5] = 0
$data[140,6] = 0
$data[140,7] = 0
$data[141,0] = 'Islas Virgenes de los Estados Unidos'
$data[141,1] = 6
$data[141,2] = 0
$data[141,3] = 0
$data[141,4] = 6
$data[141,5] = 0
$data[141,6] = 0
$data[141,7] = 0
$data[142,0] = 'San Martin (Parte Francesa)'
$data[142,1] = 5
$data[142,2] = 0
$data[142,3] = 0
$data[142,4] = 5
$data[142,5] = 0
$data[142,6] = 0
$data[142,7] = 0
$data[143,0] = 'Isla de Man'
$data[143,1] = 5
$data[143,2] = 3
$data[143,3] = 0
$data[143,4] = 5
$data[143,5] = 0
$data[143,6] = 0
$data[143,7] = 0
$data[144,0] = 'Surinam'
$data[144,1] = 5
$data[144,2] = 0
$data[144,3] = 0
$data[144,4] = 5
$data[144,5] = 0
$data[144,6] = 0
$data[144,7] = 0
$data[145,0] = 'Gabon'
$data[145,1] = 5
$data[145,2] = 0
$data[145,3] = 0
$data[145,4] = 4
$data[145,5] = 0
$data[145,6] = 0
$data[145,7] = 1
$data[146,0] = 'Nueva Caledonia'
$data[146,1] = 4
$data[146,2] = 0
$data[146,3] = 0
$data[146,4] = 4
$data[146,5] = 0
$data[146,6] = 0
$data[146,7] = 0
$data[147,0] = 'Suazilandia'
$data[147,1] = 4
$data[147,2] = 3
$data[147,3] = 0
$data[147,4] = 4
$data[147,5] = 0
$data[147,6] = 0
$data[147,7] = 0
$data[148,0] = 'Bahamas'
$data[148,1] = 4
$data[148,2] = 0
$data[148,3] = 0
$data[148,4] = 4
$data[148,5] = 0
$data[148,6] = 0
$data[148,7] = 0
$data[149,0] = 'San Bartolome'
$data[149,1] = 3
$data[149,2] = 0
$data[149,3] = 0
$data[149,4] = 3
$data[149,5] = 0
$data[149,6] = 0
$data[149,7] = 0
$data[150,0] = 'Zimbabue'
$data[150,1] = 3
$data[150,2] = 0
$data[150,3] = 0
$data[150,4] = 3
$data[150,5] = 0
$data[150,6] = 0
$data[150,7] = 0
$data[151,0] = 'El Salvador'
$data[151,1] = 3
$data[151,2] = 0
$data[151,3] = 0
$data[151,4] = 3
$data[151,5] = 0
$data[151,6] = 0
$data[151,7] = 0
$data[152,0] = 'Cabo Verde'
$data[152,1] = 3
$data[152,2] = 0
$data[152,3] = 0
$data[152,4] = 3
$data[152,5] = 0
$data[152,6] = 0
$data[152,7] = 0
$data[153,0] = 'Republica de Africa Central'
$data[153,1] = 3
$data[153,2] = 0
$data[153,3] = 0
$data[153,4] = 3
$data[153,5] = 0
$data[153,6] = 0
$data[153,7] = 0
$data[154,0] = 'Congo'
$data[154,1] = 3
$data[154,2] = 0
$data[154,3] = 0
$data[154,4] = 3
$data[154,5] = 0
$data[154,6] = 0
$data[154,7] = 0
$data[155,0] = 'Zambia'
$data[155,1] = 3
$data[155,2] = 1
$data[155,3] = 0
$data[155,4] = 3
$data[155,5] = 0
$data[155,6] = 0
$data[155,7] = 0
$data[156,0] = 'Madagascar'
$data[156,1] = 3
$data[156,2] = 0
$data[156,3] = 0
$data[156,4] = 3
$data[156,5] = 0
$data[156,6] = 0
$data[156,7] = 0
$data[157,0] = 'Namibia'
$data[157,1] = 3
$data[157,2] = 0
$data[157,3] = 0
$data[157,4] = 3
$data[157,5] = 0
$data[157,6] = 0
$data[157,7] = 0
$data[158,0] = 'Liberia'
$data[158,1] = 3
$data[158,2] = 0
$data[158,3] = 0
$data[158,4] = 3
$data[158,5] = 0
$data[158,6] = 0
$data[158,7] = 0
$data[159,0] = 'Curazao'
$data[159,1] = 3
$data[159,2] = 0
$data[159,3] = 0
$data[159,4] = 2
$data[159,5] = 0
$data[159,6] = 0
$data[159,7] = 1
$data[160,0] = 'Islas Caimanes'
$data[160,1] = 3
$data[160,2] = 0
$data[160,3] = 0
$data[160,4] = 2
$data[160,5] = 0
$data[160,6] = 0
$data[160,7] = 1
$data[161,0] = 'Bermudas'
$data[161,1] = 2
$data[161,2] = 0
$data[161,3] = 0
$data[161,4] = 2
$data[161,5] = 0
$data[161,6] = 0
$data[161,7] = 0
$data[162,0] = 'Butan'
$data[162,1] = 2
$data[162,2] = 0
$data[162,3] = 0
$data[162,4] = 2
$data[162,5] = 0
$data[162,6] = 0
$data[162,7] = 0
$data[163,0] = 'Benin'
$data[163,1] = 2
$data[163,2] = 0
$data[163,3] = 0
$data[163,4] = 2
$data[163,5] = 0
$data[163,6] = 0
$data[163,7] = 0
$data[164,0] = 'Groenlandia'
$data[164,1] = 2
$data[164,2] = 0
$data[164,3] = 0
$data[164,4] = 2
$data[164,5] = 0
$data[164,6] = 0
$data[164,7] = 0
$data[165,0] = 'Fiyi'
$data[165,1] = 2
$data[165,2] = 0
$data[165,3] = 0
$data[165,4] = 2
$data[165,5] = 0
$data[165,6] = 0
$data[165,7] = 0
$data[166,0] = 'Santa Lucia'
$data[166,1] = 2
$data[166,2] = 0
$data[166,3] = 0
$data[166,4] = 2
$data[166,5] = 0
$data[166,6] = 0
$data[166,7] = 0
$data[167,0] = 'Guinea'
$data[167,1] = 2
$data[167,2] = 0
$data[167,3] = 0
$data[167,4] = 2
$data[167,5] = 0
$data[167,6] = 0
$data[167,7] = 0
$data[168,0] = 'Nicaragua'
$data[168,1] = 2
$data[168,2] = 0
$data[168,3] = 0
$data[168,4] = 2
$data[168,5] = 0
$data[168,6] = 0
$data[168,7] = 0
$data[169,0] = 'Mauritania'
$data[169,1] = 2
$data[169,2] = 0
$data[169,3] = 0
$data[169,4] = 2
$data[169,5] = 0
$data[169,6] = 0
$data[169,7] = 0
$data[170,0] = 'Niger'
$data[170,1] = 2
$data[170,2] = 1
$data[170,3] = 0
$data[170,4] = 2
$data[170,5] = 0
$data[170,6] = 0
$data[170,7] = 0
$data[171,0] = 'Haiti'
$data[171,1] = 2
$data[171,2] = 0
$data[171,3] = 0
$data[171,4] = 2
$data[171,5] = 0
$data[171,6] = 0
$data[171,7] = 0
$data[172,0] = 'Angola'
$data[172,1] = 2
$data[172,2] = 0
$data[172,3] = 0
$data[172,4] = 2
$data[172,5] = 0
$data[172,6] = 0
$data[172,7] = 0
$data[173,0] = 'Sudan'
$data[173,1] = 2
$data[173,2] = 0
$data[173,3] = 0
$data[173,4] = 1
$data[173,5] = 0
$data[173,6] = 0
$data[173,7] = 1
$data[174,0] = 'San Vicente y las Granadinas'
$data[174,1] = 1
$data[174,2] = 0
$data[174,3] = 0
$data[174,4] = 1
$data[174,5] = 0
$data[174,6] = 0
$data[174,7] = 0
$data[175,0] = 'Timor Oriental'
$data[175,1] = 1
$data[175,2] = 0
$data[175,3] = 0
$data[175,4] = 1
$data[175,5] = 0
$data[175,6] = 0
$data[175,7] = 0
$data[176,0] = 'Republica del Chad'
$data[176,1] = 1
$data[176,2] = 0
$data[176,3] = 0
$data[176,4] = 1
$data[176,5] = 0
$data[176,6] = 0
$data[176,7] = 0
$data[177,0] = 'Eritrea'
$data[177,1] = 1
$data[177,2] = 0
$data[177,3] = 0
$data[177,4] = 1
$data[177,5] = 0
$data[177,6] = 0
$data[177,7] = 0
$data[178,0] = 'Antigua y Barbuda'
$data[178,1] = 1
$data[178,2] = 0
$data[178,3] = 0
$data[178,4] = 1
$data[178,5] = 0
$data[178,6] = 0
$data[178,7] = 0
$data[179,0] = 'Republica de Yibuti'
$data[179,1] = 1
$data[179,2] = 0
$data[179,3] = 0
$data[179,4] = 1
$data[179,5] = 0
$data[179,6] = 0
$data[179,7] = 0
$data[180,0] = 'Uganda'
$data[180,1] = 1
$data[180,2] = 0
$data[180,3] = 0
$data[180,4] = 1
$data[180,5] = 0
$data[180,6] = 0
$data[180,7] = 0
$data[181,0] = 'Papua Nueva Guinea'
$data[181,1] = 1
$data[181,2] = 0
$data[181,3] = 0
$data[181,4] = 1
$data[181,5] = 0
$data[181,6] = 0
$data[181,7] = 0
$data[182,0] = 'Mozambique'
$data[182,1] = 1
$data[182,2] = 1
$data[182,3] = 0
$data[182,4] = 1
$data[182,5] = 0
$data[182,6] = 0
$data[182,7] = 0
$data[183,0] = 'Montserrat'
$data[183,1] = 1
$data[183,2] = 0
$data[183,3] = 0
$data[183,4] = 1
$data[183,5] = 0
$data[183,6] = 0
$data[183,7] = 0
$data[184,0] = 'San Martin (Parte Holandesa)'
$data[184,1] = 1
$data[184,2] = 0
$data[184,3] = 0
$data[184,4] = 1
$data[184,5] = 0
$data[184,6] = 0
$data[184,7] = 0
$data[185,0] = 'Santa Sede'
$data[185,1] = 1
$data[185,2] = 0
$data[185,3] = 0
$data[185,4] = 1
$data[185,5] = 0
$data[185,6] = 0
$data[185,7] = 0
$data[186,0] = 'Gambia'
$data[186,1] = 1
$data[186,2] = 0
$data[186,3] = 0
$data[186,4] = 1
$data[186,5] = 0
$data[186,6] = 0
$data[186,7] = 0
$data[187,0] = 'Somalia'
$data[187,1] = 1
$data[187,2] = 0
$data[187,3] = 0
$data[187,4] = 1
$data[187,5] = 0
$data[187,6] = 0
$data[187,7] = 0
$data[188,0] = 'Nepal'
$data[188,1] = 1
$data[188,2] = 0
$data[188,3] = 1
$data[188,4] = 0
$data[188,5] = 0
$data[188,6] = 0
$data[188,7] = 0

$ws.Range("A4:H192").Value = $data